# Daily attendance processing - swap the order of names in the
# "Recorded By" column (G) from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" wherever that exact value occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
